$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2023-10-13 Friday" "2023-10-14 Saturday"

Replace-Text "29×23=667" "48×21=1008"
Replace-Text "71×31=2201" "56×25=1400"
Replace-Text "81×70=5670" "40×70=2800"
Replace-Text "91×50=4550" "95×28=2660"
Replace-Text "62×21=1302" "33×53=1749"

Replace-Text "52×62=3224" "89×30=2670"
Replace-Text "93×16=1488" "45×28=1260"
Replace-Text "20×32=640" "81×44=3564"
Replace-Text "79×66=5214" "33×19=627"
Replace-Text "26×87=2262" "85×69=5865"

Replace-Text "81×42=3402" "76×94=7144"
Replace-Text "17×80=1360" "70×58=4060"
Replace-Text "68×39=2652" "26×48=1248"
Replace-Text "78×90=7020" "23×41=943"
Replace-Text "58×21=1218" "80×66=5280"

Replace-Text "24×13=312" "84×92=7728"
Replace-Text "26×65=1690" "76×22=1672"
Replace-Text "13×20=260" "35×14=490"
Replace-Text "39×75=2925" "79×62=4898"
Replace-Text "48×75=3600" "62×56=3472"

Replace-Text "69×23=1587" "36×91=3276"
Replace-Text "15×74=1110" "83×97=8051"
Replace-Text "21×84=1764" "70×72=5040"
Replace-Text "25×44=1100" "79×74=5846"
Replace-Text "68×76=5168" "85×46=3910"
